$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Clients" to "clients"
$ws.Name = "clients"

# Clear out the "itineraire" data values in column E for rows 2-6
# (header E1 "itineraire" stays intact)
$ws.Range("E2:E6").ClearContents()
